# Leave card update - 3/9/2023
# 1) Insert a new row above the final "totals" row of Table1 (old row 130 -> new row 131)
#    so the table can grow by one row, then lay out the 2023/2024 monthly VL-earn
#    schedule in the freed-up rows (104-127), mirroring the 2022 block (rows 85-99)
#    directly above it.
# 2) Record a VL(6-4-0) leave taken 12/16-12/27/2022 (6.5 days charged) on row 103.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# ---------------------------------------------------------------------------
# Step 1: shift the special "final" row (130) down to 131 so a fresh, normally
# formatted row can take its place at 130 - mirrors Excel's native
# "insert table row above" on the last row of Table1.
# ---------------------------------------------------------------------------
$ws.Range("A130:K130").Copy()
$ws.Range("A131:K131").PasteSpecial(-4122)   # xlPasteFormats

$tbl.Resize($ws.Range("A8:K131"))
$ws.Range("G131").Formula = $ws.Range("G130").Formula

# Re-format the newly freed-up row 130 like a normal data row (copy format
# from row 129, the row directly above it).
$ws.Range("A129:K129").Copy()
$ws.Range("A130:K130").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G130").Formula = $ws.Range("G129").Formula

# ---------------------------------------------------------------------------
# Step 2: row 104 becomes the "2023" year-header row (same look as A85 = "2022")
# ---------------------------------------------------------------------------
$ws.Range("A104").Value = "'2023"
$ws.Range("A85").Copy()
$ws.Range("A104").PasteSpecial(-4122)        # xlPasteFormats

# ---------------------------------------------------------------------------
# Step 3: monthly VL-earned schedule for 2023 / 2024 (rows 105-127)
# ---------------------------------------------------------------------------
$ws.Range("A105").Value = "1/1/2023"
$ws.Range("B105").Value = "SP(1-0-00)"
$ws.Range("C105").Value = 1.25
$ws.Range("K88").Copy()
$ws.Range("K105").PasteSpecial(-4122)        # xlPasteFormats (style 49, date fmt)
$ws.Range("K105").Value = "1/3/2023"

$ws.Range("A106").Value = "2/1/2023"
$ws.Range("C106").Value = 1.25

$ws.Range("A107").Value = "3/1/2023"
$ws.Range("A108").Value = "4/1/2023"
$ws.Range("A109").Value = "5/1/2023"
$ws.Range("A110").Value = "6/1/2023"
$ws.Range("A111").Value = "7/1/2023"
$ws.Range("A112").Value = "8/1/2023"
$ws.Range("A113").Value = "9/1/2023"
$ws.Range("A114").Value = "10/1/2023"
$ws.Range("A115").Value = "11/1/2023"
$ws.Range("A116").Value = "12/1/2023"
$ws.Range("A117").Value = "1/1/2024"
$ws.Range("A118").Value = "2/1/2024"
$ws.Range("A119").Value = "3/1/2024"
$ws.Range("A120").Value = "4/1/2024"
$ws.Range("A121").Value = "5/1/2024"
$ws.Range("A122").Value = "6/1/2024"
$ws.Range("A123").Value = "7/1/2024"
$ws.Range("A124").Value = "8/1/2024"
$ws.Range("A125").Value = "9/1/2024"
$ws.Range("A126").Value = "10/1/2024"
$ws.Range("A127").Value = "11/1/2024"

# ---------------------------------------------------------------------------
# Step 4: row 103 - leave particulars / charged days / remarks
# ---------------------------------------------------------------------------
$ws.Range("B103").Value = "VL(6-4-0)"
$ws.Range("D103").Value = 6.5
$ws.Range("K103").Value = "12/16,19,20,21,22,23HD, 27"
